$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.334.25'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.800.08'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.30'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.574'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.02'
$ws.Range('E8').Value = '  +10.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.300'
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0692'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.061.78'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.70'
$ws.Range('E13').Value = '  +6.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.802.87'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.642'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.51'
$ws.Range('E16').Value = '  +5.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.353.72'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.08'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.17'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0794'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.67'
$ws.Range('E21').Value = '  +4.21%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.16'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.35'
$ws.Range('E25').Value = '  +3.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.95'
$ws.Range('E26').Value = '  +9.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.83'
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.99'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0530'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.83'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.397.69'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.670'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.47'
$ws.Range('E37').Value = '  -4.19%  '
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.24'
$ws.Range('E40').Value = '  +10.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.961'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '82.44'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.82'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.40'
$ws.Range('E45').Value = '  -3.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0509'
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.962.43'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.38'
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  +0.81%  '
